$d = $word.ActiveDocument

$d.Content.Find.Execute("769×3=", $true, $false, $false, $false, $false, $true, 1, $false, "912×5=", 2) | Out-Null
$d.Content.Find.Execute("916×4=", $true, $false, $false, $false, $false, $true, 1, $false, "514×3=", 2) | Out-Null
$d.Content.Find.Execute("443×7=", $true, $false, $false, $false, $false, $true, 1, $false, "226×6=", 2) | Out-Null
$d.Content.Find.Execute("975×4=", $true, $false, $false, $false, $false, $true, 1, $false, "357×7=", 2) | Out-Null
$d.Content.Find.Execute("983×5=", $true, $false, $false, $false, $false, $true, 1, $false, "298×6=", 2) | Out-Null
$d.Content.Find.Execute("585×5=", $true, $false, $false, $false, $false, $true, 1, $false, "752×5=", 2) | Out-Null
$d.Content.Find.Execute("792×2=", $true, $false, $false, $false, $false, $true, 1, $false, "147×4=", 2) | Out-Null
$d.Content.Find.Execute("912×3=", $true, $false, $false, $false, $false, $true, 1, $false, "179×6=", 2) | Out-Null
$d.Content.Find.Execute("832×4=", $true, $false, $false, $false, $false, $true, 1, $false, "765×4=", 2) | Out-Null
$d.Content.Find.Execute("163×2=", $true, $false, $false, $false, $false, $true, 1, $false, "428×3=", 2) | Out-Null
$d.Content.Find.Execute("926×5=", $true, $false, $false, $false, $false, $true, 1, $false, "234×2=", 2) | Out-Null
$d.Content.Find.Execute("305×3=", $true, $false, $false, $false, $false, $true, 1, $false, "341×6=", 2) | Out-Null
$d.Content.Find.Execute("391×3=", $true, $false, $false, $false, $false, $true, 1, $false, "816×8=", 2) | Out-Null
$d.Content.Find.Execute("650×6=", $true, $false, $false, $false, $false, $true, 1, $false, "303×7=", 2) | Out-Null
$d.Content.Find.Execute("275×2=", $true, $false, $false, $false, $false, $true, 1, $false, "535×4=", 2) | Out-Null
$d.Content.Find.Execute("438×5=", $true, $false, $false, $false, $false, $true, 1, $false, "754×9=", 2) | Out-Null
$d.Content.Find.Execute("187×3=", $true, $false, $false, $false, $false, $true, 1, $false, "838×3=", 2) | Out-Null
$d.Content.Find.Execute("211×3=", $true, $false, $false, $false, $false, $true, 1, $false, "722×5=", 2) | Out-Null
$d.Content.Find.Execute("935×5=", $true, $false, $false, $false, $false, $true, 1, $false, "934×5=", 2) | Out-Null
$d.Content.Find.Execute("256×9=", $true, $false, $false, $false, $false, $true, 1, $false, "636×7=", 2) | Out-Null
$d.Content.Find.Execute("985×5=", $true, $false, $false, $false, $false, $true, 1, $false, "118×6=", 2) | Out-Null
$d.Content.Find.Execute("220×4=", $true, $false, $false, $false, $false, $true, 1, $false, "826×2=", 2) | Out-Null
$d.Content.Find.Execute("652×2=", $true, $false, $false, $false, $false, $true, 1, $false, "549×4=", 2) | Out-Null
$d.Content.Find.Execute("395×9=", $true, $false, $false, $false, $false, $true, 1, $false, "536×9=", 2) | Out-Null
$d.Content.Find.Execute("682×8=", $true, $false, $false, $false, $false, $true, 1, $false, "209×6=", 2) | Out-Null
